# Insert a new weekly record as row 139 in the "Coliflor" price table.
# This shifts the existing rows 139-220 down to 140-221 (Excel's normal
# row-insert/shift-down behavior) and populates the newly created row 139
# with the new observation's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 139, pushing rows 139:220 down to 140:221.
$ws.Rows("139:139").Insert()

# Populate the new row 139 with the new weekly data point.
$ws.Range("A139").Value = 7
$ws.Range("B139").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C139").Value = 'Ñuble'
$ws.Range("D139").Value = 44582
$ws.Range("E139").Value = 16
$ws.Range("F139").Value = 100112008
$ws.Range("G139").Value = 'Coliflor'
$ws.Range("H139").Value = 'Sin especificar'
$ws.Range("I139").Value = 'Primera'
$ws.Range("J139").Value = 300
$ws.Range("K139").Value = 850
$ws.Range("L139").Value = 900
$ws.Range("M139").Value = 875
$ws.Range("N139").Value = '$/unidad'
$ws.Range("O139").Value = 'Provincia de Diguillín'
$ws.Range("P139").Value = 875
$ws.Range("Q139").Value = 1
$ws.Range("R139").Value = 'Hortaliza'
